# Commit: "updated code block docu in folder input/pm,"
#
# The original first sheet ("20180916") is duplicated and the copy is
# inserted in front of it, renamed to "20180917". In the new copy the
# code-block comment cell (C7) gets an extra link appended, the cell is
# switched to a wrap-text style and the row is heightened to show the
# extra line. The original sheet keeps its old content/selection but
# loses the "active tab" state (and its cursor moves to C7, where the
# edit happened), while the new sheet becomes the active tab with the
# cursor left at J21.

$wb = $excel.ActiveWorkbook

# The sheet we are duplicating is the current first tab, "20180916".
$source = $wb.Worksheets.Item("20180916")

# Duplicate it, Excel inserts the copy directly before $source and makes
# it the active sheet ("20180916 (2)" at this point).
$source.Copy($source)

# The newly inserted copy is now the first worksheet. Re-resolve both
# sheets by name afterwards -- index-based handles shift when a sheet is
# inserted, but name-based lookup stays pinned to the right worksheet.
$newSheet = $wb.Worksheets.Item(1)
$newSheet.Name = "20180917"
$source = $wb.Worksheets.Item("20180916")

# --- Edit the content of the new sheet -----------------------------------
# Append the forum link to the existing comment in C7.
$newSheet.Range("C7").Value = "Debug for both targets not yet possible (e.g. gdb on x86 and Segger ICD for Cortex M4)`nhttps://www.mikrocontroller.net/topic/265600"

# Wrap the text so the added line is visible...
$newSheet.Range("C7").WrapText = $true

# ...and grow the row so the wrapped text actually fits.
$newSheet.Rows.Item(7).RowHeight = 45

# --- Restore the original sheet's cursor position -------------------------
# The source sheet is no longer the active tab; its selection moved to C7
# (where the text originally lived) before focus shifted to the new sheet.
$source.Range("C7").Select()

# Leave the cursor where the author left it on the new sheet, and make it
# the active tab (selecting a range on a sheet also activates that sheet).
$newSheet.Range("J21").Select()
